# The commit swaps the two embedded theme parts of the deck:
#   ppt/theme/theme1.xml  (bound to the slide master)  "Integral"      -> "Office Theme"
#   ppt/theme/theme2.xml  (bound to the notes master)   "Office Theme" -> "Integral"
#
# The slide master's theme is reachable (and writable) through the
# PowerPoint object model via SlideMaster.Theme.ThemeColorScheme, so we
# repoint its 12 theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) from the "Integral" palette to the stock "Office" palette -
# the same edit a user makes by picking the built-in "Office" colour
# scheme from Design > Variants > Colors (the font scheme/format scheme
# for both themes were already identical, so only the colour scheme
# needs to change).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# PowerPoint ThemeColorScheme.Colors(index) -> RGB (stored BGR-packed,
# same convention as the VBA RGB() function: r + g*256 + b*65536).
$colors.Colors(1).RGB  = 0         # dk1      000000
$colors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388   # dk2      44546A
$colors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407     # accent4  FFC000
$colors.Colors(9).RGB  = 12874308  # accent5  4472C4
$colors.Colors(10).RGB = 4697456   # accent6  70AD47
$colors.Colors(11).RGB = 12673797  # hlink    0563C1
$colors.Colors(12).RGB = 7491477   # folHlink 954F72
